$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 110, shifting existing rows 110:127 down to 111:128
$ws.Rows.Item(110).Insert()

# Fill in the values for the newly inserted row 110
$ws.Cells.Item(110, 1).Value = 5
$ws.Cells.Item(110, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(110, 3).Value = "Maule"
$ws.Cells.Item(110, 4).Value = 45211
$ws.Cells.Item(110, 5).Value = 7
$ws.Cells.Item(110, 6).Value = 100112026
$ws.Cells.Item(110, 7).Value = "Haba"
$ws.Cells.Item(110, 8).Value = "Sin especificar"
$ws.Cells.Item(110, 9).Value = "Primera"
$ws.Cells.Item(110, 10).Value = 300
$ws.Cells.Item(110, 11).Value = 9000
$ws.Cells.Item(110, 12).Value = 9000
$ws.Cells.Item(110, 13).Value = 9000
$ws.Cells.Item(110, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(110, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(110, 16).Value = 360
$ws.Cells.Item(110, 17).Value = 25
$ws.Cells.Item(110, 18).Value = "Hortaliza"

# Copy the date cell style (column D) from the row below to keep the custom date format
$ws.Cells.Item(111, 4).Copy()
$ws.Cells.Item(110, 4).PasteSpecial(-4122)
$excel.CutCopyMode = $false
